# Update configuration tracking with new item categories
# Adds Commerce and Document Designer item rows to the Tracker sheet,
# restyles the data rows, narrows column A, and adds a list data
# validation on the new "action" cell in row 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: Commerce / oraclecpqo_bmClone_2 -------------------------------
$ws.Cells.Item(6, 1).Value = "NAD"
$ws.Cells.Item(6, 2).Value = "Nilesh"
$ws.Cells.Item(6, 3).Value = "Commerce"
$ws.Cells.Item(6, 4).Value = "oraclecpqo_bmClone_2"
$ws.Cells.Item(6, 6).Value = "transaction"

# --- Row 7's G cell is registered before row 6's G cell (matches the
#     original authoring order captured in the shared-string table). ----
$ws.Cells.Item(7, 7).Value = "aPI_Save_t"
$ws.Cells.Item(6, 7).Value = "save_start_step"
$ws.Cells.Item(6, 8).Value = "action"
$ws.Cells.Item(6, 5).Value = $true

# --- Row 7: Commerce / oraclecpqo_bmClone_3 -------------------------------
$ws.Cells.Item(7, 1).Value = "NAD"
$ws.Cells.Item(7, 2).Value = "Nilesh"
$ws.Cells.Item(7, 3).Value = "Commerce"
$ws.Cells.Item(7, 4).Value = "oraclecpqo_bmClone_3"
$ws.Cells.Item(7, 6).Value = "transaction"
$ws.Cells.Item(7, 8).Value = "action"
$ws.Cells.Item(7, 5).Value = $true

# --- Row 8: NAD-10759 / Document Designer ---------------------------------
$ws.Cells.Item(8, 1).Value = "NAD-10759"
$ws.Cells.Item(8, 2).Value = "Archana"
$ws.Cells.Item(8, 3).Value = "Document Designer"
$ws.Cells.Item(8, 4).Value = "oraclecpqo_bmClone_2"
$ws.Cells.Item(8, 7).Value = "Hybrid Quote Document Design - English"
$ws.Cells.Item(8, 8).Value = "doc_designer"
$ws.Cells.Item(8, 5).Value = $false

# --- Restyle the whole data range (rows 2-8) with the smaller font -------
$ws.Range("A2:H8").Font.Size = 8

# --- Narrow column A, dropping the old auto "best fit" width -------------
$ws.Columns.Item(1).ColumnWidth = 8.26

# --- List validation on the new row's action cell -------------------------
$ws.Range("H6").Validation.Add(3, 1, 1, "=INDIRECT(#REF!)")
